$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.374.48"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.884.91"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7138"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.55"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08099"
$ws.Range("E8").Value = "  +4.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3140"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08366"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "1.870.71"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7228"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.254"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.03"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.285"
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008450"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "29.374.38"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.85"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.26"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").Value = "2.124.71"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.817"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1591"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.43"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.083"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.60"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.435"
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.354"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.216"
$ws.Range("E32").Value = "  -3.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05384"
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7541"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.699"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01886"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").Value = "1.281.50"
$ws.Range("E39").Value = "  +9.59%  "
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.583"
$ws.Range("E41").Value = "  +3.45%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8944"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.48"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "110.60"
$ws.Range("E44").Value = "  +3.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  +7.40%  "
$ws.Range("D47").Value = "2.017.12"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.503"
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4375"
$ws.Range("E51").Value = "  +1.70%  "
